$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2368.1538
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2368.1538
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7104.4614
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -7644.4614

$ws.Range("H73").Value = 2368.1538
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2368.1538
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7104.4614
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -8976.4614

$ws.Range("H96").Value = 907.5
$ws.Range("I96").Value = 550.5
$ws.Range("J96").Value = 1264.5
$ws.Range("K96").Value = 1651.5
$ws.Range("L96").Value = 3793.5
$ws.Range("M96").Value = -278.5
$ws.Range("N96").Value = -6539.5

$ws.Range("H106").Value = 828.1429
$ws.Range("I106").Value = 828.1429
$ws.Range("K106").Value = 828.1429
$ws.Range("M106").Value = -197.1429000000001

$ws.Range("H125").Value = 2750
$ws.Range("I125").Value = 2750
$ws.Range("K125").Value = 24750
$ws.Range("M125").Value = -22290

$ws.Range("H138").Value = 2870.5293
$ws.Range("J138").Value = 3138.4614
$ws.Range("L138").Value = 9415.3842
$ws.Range("N138").Value = -19695.3842


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3410.0527
$ws.Range("I45").Value = 2701.2144
$ws.Range("K45").Value = 2701.2144
$ws.Range("M45").Value = -2324.2144

$ws.Range("H61").Value = 2640.5454
$ws.Range("I61").Value = 2640.5454
$ws.Range("K61").Value = 2640.5454
$ws.Range("M61").Value = -2428.5454

$ws.Range("H101").Value = 77500
$ws.Range("J101").Value = 77500
$ws.Range("L101").Value = 77500
$ws.Range("N101").Value = -83990

$ws.Range("H125").Value = 97975
$ws.Range("J125").Value = 97975
$ws.Range("L125").Value = 97975
$ws.Range("N125").Value = -107815

$ws.Range("H132").Value = 7499.5
$ws.Range("I132").Value = 7142.2856
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 21426.8568
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -18896.8568
$ws.Range("N132").Value = -35060

$ws.Range("H136").Value = 2640.5454
$ws.Range("I136").Value = 2640.5454
$ws.Range("K136").Value = 7921.6362
$ws.Range("M136").Value = -5371.6362


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 5000
$ws.Range("I15").Value = 5000
$ws.Range("K15").Value = 5000
$ws.Range("M15").Value = -4773

$ws.Range("H80").Value = 215.6
$ws.Range("I80").Value = 173.16667
$ws.Range("K80").Value = 173.16667
$ws.Range("M80").Value = 824.8333299999999

$ws.Range("H83").Value = 215.6
$ws.Range("I83").Value = 173.16667
$ws.Range("K83").Value = 865.8333500000001
$ws.Range("M83").Value = 4126.16665

$ws.Range("H96").Value = 17999.666
$ws.Range("I96").Value = 16999.5
$ws.Range("J96").Value = 20000
$ws.Range("K96").Value = 16999.5
$ws.Range("L96").Value = 20000
$ws.Range("M96").Value = -14253.5
$ws.Range("N96").Value = -25492

$ws.Range("H107").Value = 3453.9666
$ws.Range("I107").Value = 1510.4286
$ws.Range("K107").Value = 1510.4286
$ws.Range("M107").Value = 409.5714

$ws.Range("H134").Value = 2488
$ws.Range("I134").Value = 1056.1428
$ws.Range("K134").Value = 3168.4284
$ws.Range("M134").Value = -633.4284000000002


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 15000
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -20630

$ws.Range("H58").Value = 3439.476
$ws.Range("I58").Value = 2818.389
$ws.Range("K58").Value = 2818.389
$ws.Range("M58").Value = -2615.389

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H107").Value = 441.55
$ws.Range("I107").Value = 373
$ws.Range("J107").Value = 715.75
$ws.Range("K107").Value = 373
$ws.Range("L107").Value = 715.75
$ws.Range("M107").Value = 1547
$ws.Range("N107").Value = -4555.75

$ws.Range("H122").Value = 2567.5833
$ws.Range("I122").Value = 2270
$ws.Range("J122").Value = 3162.75
$ws.Range("K122").Value = 6810
$ws.Range("L122").Value = 9488.25
$ws.Range("M122").Value = -4360
$ws.Range("N122").Value = -14388.25

$ws.Range("H134").Value = 2045
$ws.Range("I134").Value = 1871.303
$ws.Range("K134").Value = 5613.909000000001
$ws.Range("M134").Value = -3078.909000000001

$ws.Range("H136").Value = 3439.476
$ws.Range("I136").Value = 2818.389
$ws.Range("K136").Value = 8455.167000000001
$ws.Range("M136").Value = -5905.167000000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43.166668
$ws.Range("I2").Value = 42.142857
$ws.Range("J2").Value = 44.6
$ws.Range("K2").Value = 252.857142
$ws.Range("L2").Value = 267.6
$ws.Range("M2").Value = -139.857142
$ws.Range("N2").Value = -493.6

$ws.Range("H36").Value = 375
$ws.Range("I36").Value = 375
$ws.Range("K36").Value = 1125
$ws.Range("M36").Value = -956

$ws.Range("H81").Value = 3146.5
$ws.Range("I81").Value = 2490
$ws.Range("J81").Value = 3365.3333
$ws.Range("K81").Value = 7470
$ws.Range("L81").Value = 10095.9999
$ws.Range("M81").Value = -6347
$ws.Range("N81").Value = -12341.9999

$ws.Range("H84").Value = 3146.5
$ws.Range("I84").Value = 2490
$ws.Range("J84").Value = 3365.3333
$ws.Range("K84").Value = 22410
$ws.Range("L84").Value = 30287.9997
$ws.Range("M84").Value = -16794
$ws.Range("N84").Value = -41519.9997

$ws.Range("H111").Value = 719.3333
$ws.Range("I111").Value = 719.3333
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2157.9999
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 909.0001000000002
$ws.Range("N111").ClearContents()

$ws.Range("H131").Value = 1055.1428
$ws.Range("I131").Value = 981
$ws.Range("K131").Value = 2943
$ws.Range("M131").Value = 2097


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 198.5
$ws.Range("I9").Value = 198.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 198.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -28.5
$ws.Range("N9").ClearContents()

$ws.Range("H31").Value = 763.1429
$ws.Range("I31").Value = 763.1429
$ws.Range("K31").Value = 763.1429
$ws.Range("M31").Value = -471.1429000000001

$ws.Range("H37").Value = 763.1429
$ws.Range("I37").Value = 763.1429
$ws.Range("K37").Value = 763.1429
$ws.Range("M37").Value = -486.1429000000001

$ws.Range("H43").Value = 8000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H63").Value = 25999.8
$ws.Range("I63").Value = 24999.5
$ws.Range("K63").Value = 24999.5
$ws.Range("M63").Value = -24313.5

$ws.Range("H66").Value = 25999.8
$ws.Range("I66").Value = 24999.5
$ws.Range("K66").Value = 74998.5
$ws.Range("M66").Value = -71566.5

$ws.Range("H80").Value = 4166
$ws.Range("I80").Value = 2246
$ws.Range("J80").Value = 8006
$ws.Range("K80").Value = 2246
$ws.Range("L80").Value = 8006
$ws.Range("M80").Value = -1248
$ws.Range("N80").Value = -10002

$ws.Range("H83").Value = 4166
$ws.Range("I83").Value = 2246
$ws.Range("J83").Value = 8006
$ws.Range("K83").Value = 11230
$ws.Range("L83").Value = 40030
$ws.Range("M83").Value = -6238
$ws.Range("N83").Value = -50014

$ws.Range("H105").Value = 28773.666
$ws.Range("J105").Value = 28773.666
$ws.Range("L105").Value = 28773.666
$ws.Range("N105").Value = -35761.666

$ws.Range("H132").Value = 99626.75
$ws.Range("I132").Value = 165074.72
$ws.Range("J132").Value = 7999.6
$ws.Range("K132").Value = 495224.16
$ws.Range("L132").Value = 23998.8
$ws.Range("M132").Value = -492694.16
$ws.Range("N132").Value = -29058.8


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1002.1539
$ws.Range("I22").Value = 678.5
$ws.Range("K22").Value = 678.5
$ws.Range("M22").Value = -383.5

$ws.Range("H27").Value = 1002.1539
$ws.Range("I27").Value = 678.5
$ws.Range("K27").Value = 678.5
$ws.Range("M27").Value = -571.5

$ws.Range("H68").Value = 7074.25
$ws.Range("I68").Value = 2198.3333
$ws.Range("K68").Value = 2198.3333
$ws.Range("M68").Value = -1449.3333

$ws.Range("H71").Value = 7074.25
$ws.Range("I71").Value = 2198.3333
$ws.Range("K71").Value = 10991.6665
$ws.Range("M71").Value = -7247.666499999999

$ws.Range("H103").Value = 14155
$ws.Range("J103").Value = 14155
$ws.Range("L103").Value = 14155
$ws.Range("N103").Value = -16499

$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 11600
$ws.Range("I33").Value = 1400
$ws.Range("J33").Value = 15000
$ws.Range("K33").Value = 1400
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = -1150
$ws.Range("N33").Value = -15500

$ws.Range("H36").Value = 11600
$ws.Range("I36").Value = 1400
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 1400
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = -1150
$ws.Range("N36").Value = -15500

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H63").Value = 32187.25
$ws.Range("I63").Value = 31500
$ws.Range("J63").Value = 32416.334
$ws.Range("K63").Value = 31500
$ws.Range("L63").Value = 32416.334
$ws.Range("M63").Value = -30876
$ws.Range("N63").Value = -33664.334

$ws.Range("H66").Value = 32187.25
$ws.Range("I66").Value = 31500
$ws.Range("J66").Value = 32416.334
$ws.Range("K66").Value = 94500
$ws.Range("L66").Value = 97249.002
$ws.Range("M66").Value = -91380
$ws.Range("N66").Value = -103489.002

$ws.Range("H103").Value = 35707.715
$ws.Range("J103").Value = 35707.715
$ws.Range("L103").Value = 35707.715
$ws.Range("N103").Value = -38051.715

$ws.Range("H107").Value = 778
$ws.Range("I107").Value = 778
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -414
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 3957.487
$ws.Range("I132").Value = 3902.0344
$ws.Range("J132").Value = 4118.3
$ws.Range("K132").Value = 11706.1032
$ws.Range("L132").Value = 12354.9
$ws.Range("M132").Value = -9176.1032
$ws.Range("N132").Value = -17414.9

$ws.Range("H136").Value = 5050.9116
$ws.Range("I136").Value = 3407.0833
$ws.Range("J136").Value = 8996.1
$ws.Range("K136").Value = 10221.2499
$ws.Range("L136").Value = 26988.3
$ws.Range("M136").Value = -7671.249899999999
$ws.Range("N136").Value = -32088.3

